$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.864.48"
$ws.Range("E2").Value = "  -2.45%  "

$ws.Range("D3").Value = "3.564.24"
$ws.Range("E3").Value = "  -3.32%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "617.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.79%  "

$ws.Range("D7").Value = "3.561.51"
$ws.Range("E7").Value = "  -3.30%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  -2.10%  "

$ws.Range("E10").Value = "  -3.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.432"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.78%  "

$ws.Range("E13").Value = "  -3.36%  "

$ws.Range("D14").Value = "4.165.98"
$ws.Range("E14").Value = "  -3.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.21%  "

$ws.Range("D16").Value = "3.560.09"
$ws.Range("E16").Value = "  -3.61%  "

$ws.Range("D17").Value = "67.979.67"
$ws.Range("E17").Value = "  -2.27%  "

$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.51%  "

$ws.Range("E20").Value = "  -2.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "455.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.67%  "

$ws.Range("E22").Value = "  -1.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.647"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.76%  "

$ws.Range("D25").Value = "3.706.86"
$ws.Range("E25").Value = "  -3.30%  "

$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000117"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.84%  "

$ws.Range("E29").Value = "  -7.38%  "

$ws.Range("E30").Value = "  -3.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.79%  "

$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.09%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.84%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.158"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.04%  "

$ws.Range("D37").Value = "3.562.98"
$ws.Range("E37").Value = "  -3.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "178.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0886"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.79%  "

$ws.Range("E43").Value = "  -7.60%  "

$ws.Range("E44").Value = "  -6.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.894"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.90%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.46%  "

$ws.Range("E50").Value = "  -5.64%  "

$ws.Range("E51").Value = "  -4.92%  "
